# Update annotations for Ying Tang
# - Row 11, column B ("politeness_score") is converted from a text "3" to a
#   proper numeric value 3 (matching the rest of the column).
# - A new row 12 is appended, reusing the text-typed "3" score value and
#   containing a new annotation record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: change B11 from text "3" to numeric 3. Everything else stays the same.
$ws.Range("B11").Value = 3

# Row 12 (new row): Annotator / politeness_score / polite_expressions /
# sentence_purpose / issue_type / id / source_file / text
$ws.Range("A12").Value = "Ying Tang"

# politeness_score needs to be stored as text "3" (not a number) to match
# the original authoring pattern, so force a text number format before
# assigning the value, then strip the format again so the cell keeps the
# default style.
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "3"
$ws.Range("B12").ClearFormats()

$ws.Range("C12").Value = "无"
$ws.Range("D12").Value = "DIS"
$ws.Range("E12").Value = "EXP"
$ws.Range("F12").Value = "d3fb2dcb-ee08-4432-9f4b-c252dbb3433f"
$ws.Range("G12").Value = "SJ3dBGZ0Z_annotated.xlsx"
$ws.Range("H12").Value = "We evaluate our method on NLP task for two reasons: 1) they are particularly well-suited for evaluating our method (naturally large output spaces) 2) we did not dispose of the computational resources to tackle tasks from other domains such as vision (e.g. Flickr100M) which requires hundreds of GPUs for weeks."
